$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update column widths (raw OOXML width = ColumnWidth + 5/6) ---
$ws.Columns.Item(2).ColumnWidth = 7.166666666666667
$ws.Columns.Item(3).ColumnWidth = 7.166666666666667
$ws.Columns.Item(7).ColumnWidth = 7.166666666666667
$ws.Columns.Item(9).ColumnWidth = 7.166666666666667
$ws.Columns.Item(10).ColumnWidth = 7.166666666666667
$ws.Columns.Item(11).ColumnWidth = 7.166666666666667
$ws.Columns.Item(12).ColumnWidth = 7.166666666666667
$ws.Columns.Item(13).ColumnWidth = 7.166666666666667
$ws.Columns.Item(15).ColumnWidth = 7.166666666666667
$ws.Columns.Item(16).ColumnWidth = 7.166666666666667
$ws.Columns.Item(17).ColumnWidth = 7.166666666666667
$ws.Columns.Item(20).ColumnWidth = 8.166666666666666
$ws.Columns.Item(22).ColumnWidth = 7.166666666666667
$ws.Columns.Item(23).ColumnWidth = 7.166666666666667
$ws.Columns.Item(24).ColumnWidth = 7.166666666666667
$ws.Columns.Item(26).ColumnWidth = 7.166666666666667
$ws.Columns.Item(27).ColumnWidth = 7.166666666666667
$ws.Columns.Item(28).ColumnWidth = 7.166666666666667
$ws.Columns.Item(29).ColumnWidth = 7.166666666666667
$ws.Columns.Item(30).ColumnWidth = 7.166666666666667
$ws.Columns.Item(34).ColumnWidth = 7.166666666666667

# --- Update data rows 2-5 with new values (row timestamps + 33 junction flow values) ---
# Row 2
$ws.Cells.Item(2, 1).Value = 45149.50694444445
$ws.Cells.Item(2, 2).Value = 23.06
$ws.Cells.Item(2, 3).Value = 15.922
$ws.Cells.Item(2, 4).Value = 4.233
$ws.Cells.Item(2, 5).Value = 48.559
$ws.Cells.Item(2, 6).Value = 40.17
$ws.Cells.Item(2, 7).Value = 18.148
$ws.Cells.Item(2, 8).Value = 60.258
$ws.Cells.Item(2, 9).Value = 27.923
$ws.Cells.Item(2, 10).Value = 11.884
$ws.Cells.Item(2, 11).Value = 18.312
$ws.Cells.Item(2, 12).Value = 19.215
$ws.Cells.Item(2, 13).Value = 20.154
$ws.Cells.Item(2, 14).Value = 5.794
$ws.Cells.Item(2, 15).Value = 18.046
$ws.Cells.Item(2, 16).Value = 25.401
$ws.Cells.Item(2, 17).Value = 15.061
$ws.Cells.Item(2, 18).Value = 3.797
$ws.Cells.Item(2, 19).Value = 2.462
$ws.Cells.Item(2, 20).Value = 267.342
$ws.Cells.Item(2, 21).Value = 50.257
$ws.Cells.Item(2, 22).Value = 16.657
$ws.Cells.Item(2, 23).Value = 33.391
$ws.Cells.Item(2, 24).Value = 17.393
$ws.Cells.Item(2, 25).Value = 2.2
$ws.Cells.Item(2, 26).Value = 29.943
$ws.Cells.Item(2, 27).Value = 14.713
$ws.Cells.Item(2, 28).Value = 13.198
$ws.Cells.Item(2, 29).Value = 15.441
$ws.Cells.Item(2, 30).Value = 19.895
$ws.Cells.Item(2, 31).Value = 3.64
$ws.Cells.Item(2, 32).Value = 53.246
$ws.Cells.Item(2, 33).Value = 9.286
$ws.Cells.Item(2, 34).Value = 20.825

# Row 3
$ws.Cells.Item(3, 1).Value = 45149.51388888889
$ws.Cells.Item(3, 2).Value = 12.011
$ws.Cells.Item(3, 3).Value = 8.34
$ws.Cells.Item(3, 4).Value = 1.669
$ws.Cells.Item(3, 5).Value = 25.493
$ws.Cells.Item(3, 6).Value = 21.11
$ws.Cells.Item(3, 7).Value = 9.452
$ws.Cells.Item(3, 8).Value = 39.208
$ws.Cells.Item(3, 9).Value = 14.543
$ws.Cells.Item(3, 10).Value = 6.226
$ws.Cells.Item(3, 11).Value = 9.461
$ws.Cells.Item(3, 12).Value = 10.235
$ws.Cells.Item(3, 13).Value = 10.637
$ws.Cells.Item(3, 14).Value = 3.021
$ws.Cells.Item(3, 15).Value = 9.399
$ws.Cells.Item(3, 16).Value = 13.231
$ws.Cells.Item(3, 17).Value = 8.089
$ws.Cells.Item(3, 18).Value = 1.591
$ws.Cells.Item(3, 19).Value = 0.926
$ws.Cells.Item(3, 20).Value = 135.754
$ws.Cells.Item(3, 21).Value = 26.395
$ws.Cells.Item(3, 22).Value = 8.676
$ws.Cells.Item(3, 23).Value = 17.438
$ws.Cells.Item(3, 24).Value = 9.336
$ws.Cells.Item(3, 25).Value = 1.1
$ws.Cells.Item(3, 26).Value = 18.485
$ws.Cells.Item(3, 27).Value = 7.663
$ws.Cells.Item(3, 28).Value = 6.982
$ws.Cells.Item(3, 29).Value = 8.162
$ws.Cells.Item(3, 30).Value = 10.635
$ws.Cells.Item(3, 31).Value = 1.294
$ws.Cells.Item(3, 32).Value = 35.483
$ws.Cells.Item(3, 33).Value = 4.786
$ws.Cells.Item(3, 34).Value = 10.847

# Row 4
$ws.Cells.Item(4, 1).Value = 45149.52083333334
$ws.Cells.Item(4, 2).Value = 13.932
$ws.Cells.Item(4, 3).Value = 10.016
$ws.Cells.Item(4, 4).Value = 1.24
$ws.Cells.Item(4, 5).Value = 29.895
$ws.Cells.Item(4, 6).Value = 24.744
$ws.Cells.Item(4, 7).Value = 10.964
$ws.Cells.Item(4, 8).Value = 42.703
$ws.Cells.Item(4, 9).Value = 16.87
$ws.Cells.Item(4, 10).Value = 7.365
$ws.Cells.Item(4, 11).Value = 11.088
$ws.Cells.Item(4, 12).Value = 12.059
$ws.Cells.Item(4, 13).Value = 12.575
$ws.Cells.Item(4, 14).Value = 3.502
$ws.Cells.Item(4, 15).Value = 10.903
$ws.Cells.Item(4, 16).Value = 15.423
$ws.Cells.Item(4, 17).Value = 9.3
$ws.Cells.Item(4, 18).Value = 1.065
$ws.Cells.Item(4, 19).Value = 0.735
$ws.Cells.Item(4, 20).Value = 158.625
$ws.Cells.Item(4, 21).Value = 30.485
$ws.Cells.Item(4, 22).Value = 10.064
$ws.Cells.Item(4, 23).Value = 20.32
$ws.Cells.Item(4, 24).Value = 10.888
$ws.Cells.Item(4, 25).Value = 1.319
$ws.Cells.Item(4, 26).Value = 20.464
$ws.Cells.Item(4, 27).Value = 8.889
$ws.Cells.Item(4, 28).Value = 7.988
$ws.Cells.Item(4, 29).Value = 9.365
$ws.Cells.Item(4, 30).Value = 12.592
$ws.Cells.Item(4, 31).Value = 0.784
$ws.Cells.Item(4, 32).Value = 38.465
$ws.Cells.Item(4, 33).Value = 5.61
$ws.Cells.Item(4, 34).Value = 12.582

# Row 5
$ws.Cells.Item(5, 1).Value = 45149.52777777778
$ws.Cells.Item(5, 2).Value = 6.25
$ws.Cells.Item(5, 3).Value = 4.36
$ws.Cells.Item(5, 4).Value = 0.77
$ws.Cells.Item(5, 5).Value = 13.28
$ws.Cells.Item(5, 6).Value = 11
$ws.Cells.Item(5, 7).Value = 4.92
$ws.Cells.Item(5, 8).Value = 22.52
$ws.Cells.Item(5, 9).Value = 7.56
$ws.Cells.Item(5, 10).Value = 3.26
$ws.Cells.Item(5, 11).Value = 4.88
$ws.Cells.Item(5, 12).Value = 5.41
$ws.Cells.Item(5, 13).Value = 5.55
$ws.Cells.Item(5, 14).Value = 1.57
$ws.Cells.Item(5, 15).Value = 4.89
$ws.Cells.Item(5, 16).Value = 6.9
$ws.Cells.Item(5, 17).Value = 4.28
$ws.Cells.Item(5, 18).Value = 0.75
$ws.Cells.Item(5, 19).Value = 0.41
$ws.Cells.Item(5, 20).Value = 67.09
$ws.Cells.Item(5, 21).Value = 13.83
$ws.Cells.Item(5, 22).Value = 4.51
$ws.Cells.Item(5, 23).Value = 9.13
$ws.Cells.Item(5, 24).Value = 4.94
$ws.Cells.Item(5, 25).Value = 0.55
$ws.Cells.Item(5, 26).Value = 10.41
$ws.Cells.Item(5, 27).Value = 3.98
$ws.Cells.Item(5, 28).Value = 3.65
$ws.Cells.Item(5, 29).Value = 4.26
$ws.Cells.Item(5, 30).Value = 5.61
$ws.Cells.Item(5, 31).Value = 0.56
$ws.Cells.Item(5, 32).Value = 20.49
$ws.Cells.Item(5, 33).Value = 2.46
$ws.Cells.Item(5, 34).Value = 5.64

# --- Delete row 6 (dataset now has one fewer sample row) ---
$ws.Rows.Item(6).Delete()
